# Update countries & provincias Spain
# - Swap the displayed country name for three adjacent row pairs
#   (the underlying per-row statistics stay with their row; only the
#   country label that had been mis-ordered moves).
# - Refresh the "Datos actualizados" timestamp string.
# - Refresh the daily COVID figures (columns B-H) for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country-name swaps -----------------------------------------------
$ws.Range("A134").Value = "Mali"
$ws.Range("A135").Value = "Angola"

$ws.Range("A139").Value = "Aruba"
$ws.Range("A140").Value = "Jordania"

$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# --- Timestamp ----------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 22:26"

# --- Row 4 : Estados Unidos ---------------------------------------------
$ws.Range("B4").Value = 6320562
$ws.Range("C4").Value = 29825
$ws.Range("D4").Value = 3558105
$ws.Range("E4").Value = 2571829
$ws.Range("G4").Value = 664
$ws.Range("H4").Value = 190628

# --- Row 10 : Sudafrica ---------------------------------------------
$ws.Range("B10").Value = 633015
$ws.Range("C10").Value = 2420
$ws.Range("D10").Value = 554887
$ws.Range("E10").Value = 63565
$ws.Range("G10").Value = 174
$ws.Range("H10").Value = 14563

# --- Row 19 : Francia ---------------------------------------------
$ws.Range("D19").Value = 87206
$ws.Range("E19").Value = 182269

# --- Row 68 : Kenia ---------------------------------------------
$ws.Range("B68").Value = 34705
$ws.Range("C68").Value = 212
$ws.Range("D68").Value = 20644
$ws.Range("E68").Value = 13476
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 585

# --- Row 95 : Guinea ---------------------------------------------
$ws.Range("B95").Value = 9579
$ws.Range("C95").Value = 53
$ws.Range("D95").Value = 8726
$ws.Range("E95").Value = 792
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 61

# --- Row 123 : Cabo Verde ---------------------------------------------
$ws.Range("B123").Value = 4125
$ws.Range("C123").Value = 77
$ws.Range("D123").Value = 3502
$ws.Range("E123").Value = 582

# --- Row 134 : now "Mali" (was Angola) ---------------------------------------------
$ws.Range("B134").Value = 2807
$ws.Range("C134").Value = 5
$ws.Range("D134").Value = 2203
$ws.Range("E134").Value = 478
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 126

# --- Row 135 : now "Angola" (was Mali) ---------------------------------------------
$ws.Range("B135").Value = 2805
$ws.Range("C135").Value = 28
$ws.Range("D135").Value = 1144
$ws.Range("E135").Value = 1548
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 113

# --- Row 139 : now "Aruba" (was Jordania) ---------------------------------------------
$ws.Range("B139").Value = 2292
$ws.Range("C139").Value = 81
$ws.Range("D139").Value = 1031
$ws.Range("E139").Value = 1248
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 13

# --- Row 140 : now "Jordania" (was Aruba) ---------------------------------------------
$ws.Range("B140").Value = 2233
$ws.Range("C140").Value = 72
$ws.Range("D140").Value = 1648
$ws.Range("E140").Value = 570
$ws.Range("H140").Value = 15

# --- Row 144 : Sierra Leona ---------------------------------------------
$ws.Range("B144").Value = 2035
$ws.Range("C144").Value = 6
$ws.Range("E144").Value = 364

# --- Row 153 : Republica de Chipre ---------------------------------------------
$ws.Range("B153").Value = 1498
$ws.Range("C153").Value = 3
$ws.Range("E153").Value = 338

# --- Row 159 : Liberia ---------------------------------------------
$ws.Range("B159").Value = 1306
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 1163

# --- Row 192 : Monaco ---------------------------------------------
$ws.Range("B192").Value = 143
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 94
$ws.Range("E192").Value = 48

# --- Row 214 : now "Montserrat" (was Islas Malvinas) ---------------------------------------------
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Row 215 : now "Islas Malvinas" (was Montserrat) ---------------------------------------------
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
